$wb = $excel.ActiveWorkbook

# The invoice sheet "2025-03-14 - 25-24831" becomes "2025-05-06 - 25-24927"
$ws = $wb.Worksheets.Item("2025-03-14 - 25-24831")

# Rename the sheet tab (this also updates the TitlesOfParts entries
# automatically, but the _xlnm.Print_Area defined name text needs a nudge
# below to pick up the new sheet name).
$ws.Name = "2025-05-06 - 25-24927"

# Re-assert the print area so the _xlnm.Print_Area defined name is rewritten
# with the new sheet name (it otherwise keeps referencing the old name).
$ws.PageSetup.PrintArea = '$A$1:$F$88'

# Update the invoice header texts
$ws.Range("B21").Value = "Le 6 MAI 2025"
$ws.Range("E28").Value = "25-24927"
$ws.Range("B34").Formula = "' - TEST"

# Update the invoice line amounts
$ws.Range("C66").Value = 9.9
$ws.Range("E69").Value = 3811.5
$ws.Range("E73").Value = 3811.5
$ws.Range("E74").Value = 190.58
$ws.Range("E75").Value = 380.2
$ws.Range("E77").Value = 4382.28
$ws.Range("E81").Value = 4382.28
